$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = 6329
$ws.Range("C23").Value = 996
$ws.Range("D23").Value = 5891358
$ws.Range("E23").Value = 930.8513193237478
$ws.Range("F23").Value = 8.596431022649288
$ws.Range("G23").Value = 3.642039542143594
$ws.Range("H23").Value = 26.23672221873432
